$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.570.17'
$ws.Range('E2').Value = '  -7.13%  '
$ws.Range('D3').Value = '1.690.87'
$ws.Range('E3').Value = '  -5.66%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.99'
$ws.Range('E5').Value = '  -5.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5110'
$ws.Range('E6').Value = '  -13.11%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2664'
$ws.Range('E8').Value = '  -3.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '22.05'
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06307'
$ws.Range('E10').Value = '  -6.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07354'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').Value = '1.693.73'
$ws.Range('E12').Value = '  -5.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.518'
$ws.Range('E13').Value = '  -5.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5774'
$ws.Range('E14').Value = '  -5.80%  '
$ws.Range('D15').Value = '1.919.72'
$ws.Range('E15').Value = '  -5.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008540'
$ws.Range('E16').Value = '  -4.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.25'
$ws.Range('E17').Value = '  -13.41%  '
$ws.Range('D18').Value = '26.599.02'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.990'
$ws.Range('E19').Value = '  -7.87%  '
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('E21').Value = '  -4.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '186.62'
$ws.Range('E22').Value = '  -10.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.258'
$ws.Range('E23').Value = '  -8.40%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.76'
$ws.Range('E25').Value = '  -5.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.496'
$ws.Range('E26').Value = '  -7.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1170'
$ws.Range('E27').Value = '  -7.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.79'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.341'
$ws.Range('E29').Value = '  -5.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05733'
$ws.Range('E30').Value = '  -8.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.341'
$ws.Range('E31').Value = '  -5.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.520'
$ws.Range('E32').Value = '  -6.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.511'
$ws.Range('E33').Value = '  -7.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.643'
$ws.Range('E34').Value = '  -5.43%  '
$ws.Range('E35').Value = '  -2.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5984'
$ws.Range('E36').Value = '  -6.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.367'
$ws.Range('E37').Value = '  -5.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.675'
$ws.Range('E38').Value = '  -1.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01620'
$ws.Range('E39').Value = '  -4.34%  '
$ws.Range('D40').Value = '1.090.95'
$ws.Range('E40').Value = '  -4.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8601'
$ws.Range('E41').Value = '  -1.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.833'
$ws.Range('E42').Value = '  -8.68%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.57'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').Value = '1.846.33'
$ws.Range('E45').Value = '  -5.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000117'
$ws.Range('E46').Value = '  +5.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.33'
$ws.Range('E47').Value = '  -5.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.073'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4323'
$ws.Range('E50').Value = '  -3.64%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05231'
$ws.Range('E51').Value = '  -4.32%  '
